# Updated data for care, v57 (added Viattence)
#
# Insert a new organisation row ("Viattence") into the routekaart status
# list, in its correct alphabetical position between "Verpleeghuis
# Bergweide (Stichting)" (row 222) and "Viersprong ..." (old row 223,
# now row 224), with status "Actueel en vastgesteld" and styled with a
# 9pt Verdana font on the organisation-name cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 223 (and everything below it) down by one to make room.
$ws.Rows("223:223").Insert()

# Populate the new row.
$ws.Range("A223").Value = "Viattence"
$ws.Range("B223").Value = "Actueel en vastgesteld"

# Style the organisation-name cell with the smaller Verdana font used for
# this entry.
$ws.Range("A223").Font.Name = "Verdana"
$ws.Range("A223").Font.Size = 9

# Leave the selection on the newly added cell, matching the author's
# working state when the change was saved.
$ws.Range("A223").Select()
